# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# This script updates the "resendAll.230511" API template workbook so that:
#  - The "Body" (Request) sheet's single body row now documents the request
#    as a schema reference ("resendAll.230511Request") instead of the old
#    ad-hoc "dateTime" field, and the now-superseded extra body rows are
#    removed.
#  - The "200" (Response) sheet is likewise turned into a schema reference
#    ("resendAll.230511Response") and its extra content rows are removed.
#  - The "204" sheet gains a row documenting the same response schema
#    ("resendAll.230511Response").
#  - The "400" sheet becomes a schema reference to the shared "errorResponse"
#    schema, and its extra content rows are removed.
#  - The "401", "403", "404", "429" and "500" sheets each gain a row
#    documenting the shared "errorResponse1" schema.

$wb = $excel.ActiveWorkbook

function Set-SchemaRow3 {
    param(
        $ws,
        [string]$SectionValue,
        [string]$SchemaName
    )

    $ws.Range("A3").Value = $SectionValue      # Section
    $ws.Range("B3").Value = $SchemaName        # Name
    $ws.Range("D3").ClearContents()            # Description
    $ws.Range("E3").Value = "schema"           # Type
    $ws.Range("G3").Value = $SchemaName        # Schema Name
    $ws.Range("I3").Value = "Yes"              # Mandatory
    $ws.Range("L3").ClearContents()            # PatternEba
    $ws.Range("O3").ClearContents()            # Example
}

# ---------------------------------------------------------------------------
# Body (Request) sheet: row 3 becomes the request-schema reference, rows
# 4-8 (fileType / fileStatus / lacNumber / dateFrom / dateTo) are removed.
# ---------------------------------------------------------------------------
$wsBody = $wb.Worksheets.Item("Body")
Set-SchemaRow3 $wsBody "body" "resendAll.230511Request"
$wsBody.Range("A4:O8").ClearContents()

# ---------------------------------------------------------------------------
# 200 sheet: row 3 becomes the response-schema reference, rows 4-5
# (commandRef / commandStatus) are removed.
# ---------------------------------------------------------------------------
$ws200 = $wb.Worksheets.Item("200")
Set-SchemaRow3 $ws200 "content" "resendAll.230511Response"
$ws200.Range("A4:O5").ClearContents()

# ---------------------------------------------------------------------------
# 204 sheet: gains a new row 3 referencing the same response schema.
# ---------------------------------------------------------------------------
$ws204 = $wb.Worksheets.Item("204")
Set-SchemaRow3 $ws204 "content" "resendAll.230511Response"

# ---------------------------------------------------------------------------
# 400 sheet: row 3 becomes the shared error-schema reference, rows 4-6
# (errorCode / errorCodeDescription / requestId) are removed.
# ---------------------------------------------------------------------------
$ws400 = $wb.Worksheets.Item("400")
Set-SchemaRow3 $ws400 "content" "errorResponse"
$ws400.Range("A4:O6").ClearContents()

# ---------------------------------------------------------------------------
# 401 / 403 / 404 / 429 / 500: each gains a new row 3 referencing the
# shared "errorResponse1" schema.
# ---------------------------------------------------------------------------
foreach ($name in @("401", "403", "404", "429", "500")) {
    $ws = $wb.Worksheets.Item($name)
    Set-SchemaRow3 $ws "content" "errorResponse1"
}
